$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.130.25'
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').Value = '3.387.24'
$ws.Range('E3').Value = '  -0.64%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '''253.87'
$ws.Range('E5').Value = '  -0.51%  '

$ws.Range('D6').Value = '''662.35'
$ws.Range('E6').Value = '  -0.18%  '

$ws.Range('D7').Value = '''1.47'
$ws.Range('E7').Value = '  +1.39%  '

$ws.Range('E8').Value = '  -1.92%  '

$ws.Range('E9').Value = '  -0.49%  '

$ws.Range('E10').Value = '  +0.00%  '

$ws.Range('D11').Value = '3.383.12'
$ws.Range('E11').Value = '  -0.65%  '

$ws.Range('E12').Value = '  -2.50%  '

$ws.Range('D13').Value = '''41.89'
$ws.Range('E13').Value = '  -1.33%  '

$ws.Range('D14').Value = '97.905.28'
$ws.Range('E14').Value = '  +0.10%  '

$ws.Range('D15').Value = '''6.13'
$ws.Range('E15').Value = '  -6.02%  '

$ws.Range('E16').Value = '  -3.75%  '

$ws.Range('D17').Value = '4.014.19'
$ws.Range('E17').Value = '  -0.64%  '

$ws.Range('D18').Value = '''8.97'
$ws.Range('E18').Value = '  -0.19%  '

$ws.Range('D19').Value = '3.377.46'
$ws.Range('E19').Value = '  -0.85%  '

$ws.Range('D20').Value = '''18.13'
$ws.Range('E20').Value = '  +2.37%  '

$ws.Range('D21').Value = '''0.531'
$ws.Range('E21').Value = '  -5.15%  '

$ws.Range('D22').Value = '''10.99'
$ws.Range('E22').Value = '  -0.06%  '

$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value = '''3.45'
$ws.Range('E23').Value = '  +0.05%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '''512.85'

$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '''0.0000202'
$ws.Range('E25').Value = '  -2.11%  '

$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = '''6.95'
$ws.Range('E26').Value = '  +5.48%  '

$ws.Range('D27').Value = '''96.80'
$ws.Range('E27').Value = '  -3.35%  '

$ws.Range('D28').Value = '''12.38'
$ws.Range('E28').Value = '  -3.93%  '

$ws.Range('D29').Value = '3.566.34'

$ws.Range('D30').Value = '''11.45'
$ws.Range('E30').Value = '  -0.87%  '

$ws.Range('D31').Value = '''0.144'
$ws.Range('E31').Value = '  -3.35%  '

$ws.Range('E33').Value = '  -5.62%  '

$ws.Range('E34').Value = '  +9.06%  '

$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  -0.11%  '

$ws.Range('D36').Value = '''0.563'
$ws.Range('E36').Value = '  -2.53%  '

$ws.Range('D37').Value = '''29.01'
$ws.Range('E37').Value = '  -2.80%  '

$ws.Range('E38').Value = '  -0.15%  '

$ws.Range('E39').Value = '  -1.42%  '

$ws.Range('D40').Value = '''535.34'
$ws.Range('E40').Value = '  -0.13%  '

$ws.Range('E41').Value = '  +0.01%  '

$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('E43').Value = '  -1.20%  '

$ws.Range('D44').Value = '''0.856'
$ws.Range('E44').Value = '  -2.70%  '

$ws.Range('D45').Value = '''0.0432'
$ws.Range('E45').Value = '  +1.25%  '

$ws.Range('D46').Value = '''1.74'
$ws.Range('E46').Value = '  +0.24%  '

$ws.Range('E47').Value = '  +0.47%  '

$ws.Range('D48').Value = '''2.25'
$ws.Range('E48').Value = '  +6.52%  '

$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = '''5.63'
$ws.Range('E49').Value = '  -3.95%  '

$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '''56.18'
$ws.Range('E50').Value = '  +3.91%  '

$ws.Range('D51').Value = '''8.63'
$ws.Range('E51').Value = '  -5.76%  '
